$d = $word.ActiveDocument

$replacements = @(
    @{old = "46×47="; new = "15×19="},
    @{old = "46×91="; new = "81×82="},
    @{old = "35×12="; new = "34×13="},
    @{old = "23×27="; new = "53×27="},
    @{old = "19×98="; new = "41×77="},
    @{old = "57×64="; new = "63×99="},
    @{old = "23×65="; new = "87×80="},
    @{old = "73×37="; new = "27×78="},
    @{old = "49×41="; new = "91×13="},
    @{old = "92×30="; new = "62×52="},
    @{old = "31×40="; new = "87×50="},
    @{old = "26×52="; new = "61×16="},
    @{old = "74×23="; new = "62×73="},
    @{old = "73×88="; new = "35×58="},
    @{old = "95×33="; new = "35×65="},
    @{old = "37×11="; new = "94×76="},
    @{old = "88×91="; new = "82×28="},
    @{old = "23×30="; new = "61×96="},
    @{old = "31×83="; new = "32×62="},
    @{old = "29×87="; new = "52×26="},
    @{old = "98×98="; new = "72×52="},
    @{old = "46×80="; new = "37×92="},
    @{old = "63×78="; new = "68×84="},
    @{old = "21×74="; new = "16×44="},
    @{old = "15×28="; new = "91×91="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
